$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Kavu Furens', ['{1}{R}{G}', 'Animal — Kavu', 'Festinatio (animali oppugnare et {T} licet ordine penes tuum venit.)', 'Kavu Furentem adhibere tibi licet quandocumque subitum adhibere potes.', '3/1'])"

$ws.Range("A3:A7").EntireRow.Delete()
